$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Castillo - Cafeplaza (Red Sun), row 9: new palette offset + mark as edited
$ws.Range("D9").Value = "0x523544"
$ws.Range("F9").Value = "Sim"

# Castillo - Cafeplaza (Blue Moon), row 20: new palette offset + mark as edited
$ws.Range("D20").Value = "0x5231E8"
$ws.Range("F20").Value = "Sim"

# EletroVilla - Placa perto da estacao, also marked as edited in both versions
$ws.Range("F6").Value = "Sim"
$ws.Range("F17").Value = "Sim"
